# Work in progress on geocoding functionality
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Fill in row 9 with the new "Filtering" requirement, following the same
# pattern used by the rows above it (e.g. row 8).
$ws.Range("A9").Value = 6
$ws.Range("B9").Value = "Filtering"
$ws.Range("C9").Value = "Filtering the results by category, region, etc"
$ws.Range("D9").Value = "high"
$ws.Range("F9").Value = "4h"

# Update the active selection to F9, matching the saved selection state.
$ws.Range("F9").Select()
